# Remove ALL horizontal lines from the document.
#
# Every horizontal rule in this document is represented as its own,
# otherwise-empty paragraph holding a single legacy VML "hr" shape:
#   <w:p><w:r><w:pict><v:rect .../></w:pict></w:r></w:p>
# These paragraphs carry no visible text -- Range.Text is just the
# trailing paragraph mark -- so we find every paragraph whose text is
# empty once trimmed and delete it outright (the delete removes the
# run/pict AND the paragraph mark). We walk from the last paragraph to
# the first so that not-yet-visited indices stay valid while we mutate
# the collection.
#
# One special case: when the horizontal rule immediately preceding the
# very last paragraph in the document is removed, that trailing
# paragraph (here, the "Date: ..." line) is re-styled from
# "First Paragraph" to "Body Text" as part of the same revision, while
# keeping its existing direct character formatting (italics) intact.

$d = $word.ActiveDocument

$restyleTrailingParagraph = $false

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $trimmed = $p.Range.Text.Trim()
    if ($trimmed.Length -ne 0) {
        continue
    }

    # Horizontal-rule placeholder paragraph -- schedule the restyle of
    # the document's final paragraph if this rule sits directly in
    # front of it (i.e. nothing but this rule separates it from the
    # end of the document).
    if ($i -eq ($count - 1)) {
        $trailing = $d.Paragraphs.Item($count)
        if ($trailing.Style.NameLocal -eq "First Paragraph") {
            $restyleTrailingParagraph = $true
        }
    }

    $p.Range.Delete()
}

if ($restyleTrailingParagraph) {
    $newCount = $d.Paragraphs.Count
    $last = $d.Paragraphs.Item($newCount)

    $wasItalic = $last.Range.Font.Italic

    $last.Style = "Body Text"

    if ($wasItalic) {
        # Re-fetch after the structural style-change mutation and
        # reapply italics: first to the run text only (range excludes
        # the trailing paragraph mark), then once more across the
        # whole paragraph range to restore the complex-script italic
        # flag too. This two-step order avoids stamping stray direct
        # formatting onto the paragraph mark itself.
        $p1 = $d.Paragraphs.Item($newCount)
        $s = $p1.Range.Start
        $e = $p1.Range.End
        $d.Range($s, $e - 1).Font.Italic = $true

        $p2 = $d.Paragraphs.Item($newCount)
        $p2.Range.Font.ItalicBi = $true
    }
}
